$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$count = 0
for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value2
    if ($v -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
        $count = $count + 1
    }
}
Write-Host "Changed:" $count
